# texts.xlsx - "Translation" sheet update.
#
# The old row 58 (SingleUseId61 / PadNumeric / Left / LTR / "0") is removed,
# which shifts rows 59-62 up to become rows 58-61. A brand new row is then
# appended as row 62 for the newly added single-use string
# "SingleUseId67" (PadNumeric / Left / LTR / "60") - the numeric keypad
# default value used by the new pump-drive-controller task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the obsolete row (everything below shifts up by one).
$ws.Rows(58).Delete()

# Append the new translation row at the now-empty row 62.
$ws.Range("B62").Value = "SingleUseId67"
$ws.Range("C62").Value = "PadNumeric"
$ws.Range("D62").Value = "Left"
$ws.Range("E62").Value = "LTR"

# Force the numeric-looking text to stay a text/shared-string cell (matches
# every other "0"/"60" style value in this column) instead of being
# auto-coerced to a number, then drop back to the default style so no
# stray per-cell formatting is left behind.
$ws.Range("F62").NumberFormat = "@"
$ws.Range("F62").Value = "60"
$ws.Range("F62").Style = "Normal"
